$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '26.211.17'
Set-TextValue $ws.Range("E2") '  -3.95%  '
Set-TextValue $ws.Range("D3") '1.659.51'
Set-TextValue $ws.Range("E3") '  -2.67%  '
Set-TextValue $ws.Range("E4") '  +0.26%  '
Set-TextValue $ws.Range("D5") '217.68'
Set-TextValue $ws.Range("E5") '  -2.62%  '
Set-TextValue $ws.Range("D6") '0.5156'
Set-TextValue $ws.Range("E6") '  -2.99%  '
Set-TextValue $ws.Range("D7") '1.006'
Set-TextValue $ws.Range("E7") '  +0.34%  '
Set-TextValue $ws.Range("D8") '0.2578'
Set-TextValue $ws.Range("E8") '  -3.10%  '
Set-TextValue $ws.Range("D9") '0.06449'
Set-TextValue $ws.Range("E9") '  -2.24%  '
Set-TextValue $ws.Range("E10") '  -3.81%  '
Set-TextValue $ws.Range("D11") '0.07825'
Set-TextValue $ws.Range("E11") '  +2.53%  '
Set-TextValue $ws.Range("D12") '1.661.18'
Set-TextValue $ws.Range("E12") '  -2.52%  '
Set-TextValue $ws.Range("D13") '4.299'
Set-TextValue $ws.Range("E13") '  -4.55%  '
Set-TextValue $ws.Range("D14") '1.888.40'
Set-TextValue $ws.Range("E14") '  -2.57%  '
Set-TextValue $ws.Range("D15") '0.5547'
Set-TextValue $ws.Range("E15") '  -4.53%  '
Set-TextValue $ws.Range("D16") '0.0₅8062'
Set-TextValue $ws.Range("E16") '  -1.26%  '
Set-TextValue $ws.Range("D17") '64.33'
Set-TextValue $ws.Range("E17") '  -4.65%  '
Set-TextValue $ws.Range("D18") '26.238.49'
Set-TextValue $ws.Range("E18") '  -3.82%  '
Set-TextValue $ws.Range("B19") 'Dai'
Set-TextValue $ws.Range("C19") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D19") '1.005'
Set-TextValue $ws.Range("E19") '  +0.24%  '
Set-TextValue $ws.Range("B20") 'BitcoinCash'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D20") '211.17'
Set-TextValue $ws.Range("E20") '  -2.12%  '
Set-TextValue $ws.Range("D21") '4.432'
Set-TextValue $ws.Range("E21") '  -4.26%  '
Set-TextValue $ws.Range("D22") '10.09'
Set-TextValue $ws.Range("E22") '  -2.58%  '
Set-TextValue $ws.Range("D23") '6.035'
Set-TextValue $ws.Range("E23") '  +1.07%  '
Set-TextValue $ws.Range("E24") '  +0.21%  '
Set-TextValue $ws.Range("D25") '145.01'
Set-TextValue $ws.Range("E25") '  +0.68%  '
Set-TextValue $ws.Range("D26") '1.755'
Set-TextValue $ws.Range("E26") '  +3.09%  '
Set-TextValue $ws.Range("D27") '0.1172'
Set-TextValue $ws.Range("E27") '  -2.34%  '
Set-TextValue $ws.Range("D28") '6.998'
Set-TextValue $ws.Range("E28") '  -2.95%  '
Set-TextValue $ws.Range("D29") '15.78'
Set-TextValue $ws.Range("E29") '  -2.35%  '
Set-TextValue $ws.Range("D30") '0.05213'
Set-TextValue $ws.Range("E30") '  -3.00%  '
Set-TextValue $ws.Range("E31") '  -2.50%  '
Set-TextValue $ws.Range("D32") '3.352'
Set-TextValue $ws.Range("E32") '  -3.28%  '
Set-TextValue $ws.Range("D33") '3.231'
Set-TextValue $ws.Range("E33") '  -5.02%  '
Set-TextValue $ws.Range("D34") '1.575'
Set-TextValue $ws.Range("E34") '  -4.13%  '
Set-TextValue $ws.Range("D35") '2.762'
Set-TextValue $ws.Range("E35") '  -3.52%  '
Set-TextValue $ws.Range("D36") '0.9328'
Set-TextValue $ws.Range("E36") '  -1.57%  '
Set-TextValue $ws.Range("D37") '2.374'
Set-TextValue $ws.Range("E37") '  -1.44%  '
Set-TextValue $ws.Range("D38") '1.175.79'
Set-TextValue $ws.Range("E38") '  +12.92%  '
Set-TextValue $ws.Range("D39") '0.5708'
Set-TextValue $ws.Range("E39") '  -2.32%  '
Set-TextValue $ws.Range("D40") '0.01596'
Set-TextValue $ws.Range("E40") '  -2.47%  '
Set-TextValue $ws.Range("E41") '  +0.23%  '
Set-TextValue $ws.Range("D42") '0.8393'
Set-TextValue $ws.Range("E42") '  -0.17%  '
Set-TextValue $ws.Range("D43") '5.681'
Set-TextValue $ws.Range("E43") '  -2.12%  '
Set-TextValue $ws.Range("D44") '100.65'
Set-TextValue $ws.Range("D45") '1.798.41'
Set-TextValue $ws.Range("E45") '  -2.56%  '
Set-TextValue $ws.Range("E46") '  +3.02%  '
Set-TextValue $ws.Range("D47") '0.4540'
Set-TextValue $ws.Range("E47") '  +0.53%  '
Set-TextValue $ws.Range("D48") '55.99'
Set-TextValue $ws.Range("E48") '  -3.05%  '
Set-TextValue $ws.Range("D49") '1.006'
Set-TextValue $ws.Range("E49") '  +0.43%  '
Set-TextValue $ws.Range("D50") '7.912'
Set-TextValue $ws.Range("E50") '  -1.65%  '
Set-TextValue $ws.Range("D51") '0.05062'
Set-TextValue $ws.Range("E51") '  -3.14%  '
